$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 of the "Rules" sheet ("R40") gets its rule-name cell (B11) replaced
# with the text "1" (still a text/shared-string cell, not a number), while
# every other attribute of the cell (style s="23", row height, etc.) is left
# untouched.
#
# A plain  $ws.Range("B11").Value = "1"  would let Excel auto-detect the
# value as a *number*, which also changes the cell's style id (it picks up a
# number-formatted style). To keep B11 textual with its original style, we
# stage the literal text in an unused scratch cell (forced to text via
# NumberFormat "@"), copy it, and paste-special *values only* onto B11 - this
# swaps in the new text content but leaves B11's existing style/format alone.
# The scratch cell is then fully cleared (contents + formatting) so it leaves
# no trace in the saved workbook.
$scratch = $ws.Range("F11")
$scratch.NumberFormat = "@"
$scratch.Value = "1"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
$excel.CutCopyMode = 0
